$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update B17 value from "RBI (India)" to the new scenario text
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Update the active selection to B17 as reflected in the saved view state
$ws.Activate()
$ws.Range("B17").Select()
